$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: the plain-text mention of GENie (with its raw URL typed
# inline) becomes a real hyperlink, same as the neighbouring "Fips"
# mention just before it:
#
#   "or GENie(https://github.com/bkaradzic/GENie)(Lua based premake)..."
#
# becomes
#
#   "or" + " " + [hyperlink: GENie -> https://github.com/bkaradzic/GENie]
#   + "(Lua based premake)..."
# ------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute("GENie(https://github.com/bkaradzic/GENie)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $full.Start

    # Collapse the matched "GENie(https://github.com/bkaradzic/GENie)" span
    # into a hyperlink whose visible text is just "GENie".
    $h = $d.Hyperlinks.Add($full, "https://github.com/bkaradzic/GENie", `
        [Type]::Missing, [Type]::Missing, "GENie")

    # The word "or" and the following space (which sit immediately before
    # the new hyperlink) were merged into one run; split them back into two
    # separate runs ("or" and " ") to mirror the surrounding run structure.
    $orOnly = $d.Range($matchStart - 3, $matchStart - 1)
    if ($orOnly.Text -eq "or") {
        $orOnly.Bold = 1
        $orOnly.Bold = 0
    }
}

# ------------------------------------------------------------------
# Change 2: "Resumes for download" -> "Resume for download"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Resumes for download", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Resume for download", 2) | Out-Null

# ------------------------------------------------------------------
# Change 3: rename the bookmark behind that heading to match
# (resumes-for-download -> resume-for-download)
# ------------------------------------------------------------------
$oldBookmarkName = "resumes-for-download"
$newBookmarkName = "resume-for-download"
if ($d.Bookmarks.Exists($oldBookmarkName)) {
    $b = $d.Bookmarks($oldBookmarkName)
    $bRange = $b.Range
    $b.Delete()
    $d.Bookmarks.Add($newBookmarkName, $bRange) | Out-Null
}
